$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set entire used range to text format so numeric-looking strings stay as text
$ws.Range("A1:K9").NumberFormat = "@"

# Header row
$ws.Range("A1").Value = "venue"
$ws.Range("B1").Value = "date"
$ws.Range("C1").Value = "result"
$ws.Range("D1").Value = "ownTeam"
$ws.Range("E1").Value = "oppTeam"
$ws.Range("F1").Value = "batsman"
$ws.Range("G1").Value = "totalRuns"
$ws.Range("H1").Value = "totalBalls"
$ws.Range("I1").Value = "total4s"
$ws.Range("J1").Value = "total6s"
$ws.Range("K1").Value = "sr"

# Row 2
$ws.Range("A2").Value = " Sharjah"
$ws.Range("B2").Value = " September 27 2020"
$ws.Range("C2").Value = "Royals won by 4 wickets (with 3 balls remaining)"
$ws.Range("D2").Value = "Rajasthan Royals"
$ws.Range("E2").Value = "Kings XI Punjab"
$ws.Range("F2").Value = "Riyan Parag "
$ws.Range("G2").Value = "0"
$ws.Range("H2").Value = "2"
$ws.Range("I2").Value = "0"
$ws.Range("J2").Value = "0"
$ws.Range("K2").Value = "0.00"

# Row 3
$ws.Range("A3").Value = " Dubai (DSC)"
$ws.Range("B3").Value = " September 30 2020"
$ws.Range("C3").Value = "KKR won by 37 runs"
$ws.Range("D3").Value = "Rajasthan Royals"
$ws.Range("E3").Value = "Kolkata Knight Riders"
$ws.Range("F3").Value = "Riyan Parag "
$ws.Range("G3").Value = "1"
$ws.Range("H3").Value = "6"
$ws.Range("I3").Value = "0"
$ws.Range("J3").Value = "0"
$ws.Range("K3").Value = "16.66"

# Row 4
$ws.Range("A4").Value = " Abu Dhabi"
$ws.Range("B4").Value = " October 03 2020"
$ws.Range("C4").Value = "RCB won by 8 wickets (with 5 balls remaining)"
$ws.Range("D4").Value = "Rajasthan Royals"
$ws.Range("E4").Value = "Royal Challengers Bangalore"
$ws.Range("F4").Value = "Riyan Parag "
$ws.Range("G4").Value = "16"
$ws.Range("H4").Value = "18"
$ws.Range("I4").Value = "1"
$ws.Range("J4").Value = "0"
$ws.Range("K4").Value = "88.88"

# Row 5
$ws.Range("A5").Value = " Dubai (DSC)"
$ws.Range("B5").Value = " October 14 2020"
$ws.Range("C5").Value = "Capitals won by 13 runs"
$ws.Range("D5").Value = "Rajasthan Royals"
$ws.Range("E5").Value = "Delhi Capitals"
$ws.Range("F5").Value = "Riyan Parag "
$ws.Range("G5").Value = "1"
$ws.Range("H5").Value = "2"
$ws.Range("I5").Value = "0"
$ws.Range("J5").Value = "0"
$ws.Range("K5").Value = "50.00"

# Row 6
$ws.Range("A6").Value = " Dubai (DSC)"
$ws.Range("B6").Value = " October 22 2020"
$ws.Range("C6").Value = "Sunrisers won by 8 wickets (with 11 balls remaining)"
$ws.Range("D6").Value = "Rajasthan Royals"
$ws.Range("E6").Value = "Sunrisers Hyderabad"
$ws.Range("F6").Value = "Riyan Parag "
$ws.Range("G6").Value = "20"
$ws.Range("H6").Value = "12"
$ws.Range("I6").Value = "2"
$ws.Range("J6").Value = "1"
$ws.Range("K6").Value = "166.66"

# Row 7
$ws.Range("A7").Value = " Dubai (DSC)"
$ws.Range("B7").Value = " November 01 2020"
$ws.Range("C7").Value = "KKR won by 60 runs"
$ws.Range("D7").Value = "Rajasthan Royals"
$ws.Range("E7").Value = "Kolkata Knight Riders"
$ws.Range("F7").Value = "Riyan Parag "
$ws.Range("G7").Value = "0"
$ws.Range("H7").Value = "7"
$ws.Range("I7").Value = "0"
$ws.Range("J7").Value = "0"
$ws.Range("K7").Value = "0.00"

# Row 8
$ws.Range("A8").Value = " Sharjah"
$ws.Range("B8").Value = " September 22 2020"
$ws.Range("C8").Value = "Royals won by 16 runs"
$ws.Range("D8").Value = "Rajasthan Royals"
$ws.Range("E8").Value = "Chennai Super Kings"
$ws.Range("F8").Value = "Riyan Parag "
$ws.Range("G8").Value = "6"
$ws.Range("H8").Value = "4"
$ws.Range("I8").Value = "1"
$ws.Range("J8").Value = "0"
$ws.Range("K8").Value = "150.00"

# Row 9
$ws.Range("A9").Value = " Dubai (DSC)"
$ws.Range("B9").Value = " October 11 2020"
$ws.Range("C9").Value = "Royals won by 5 wickets (with 1 ball remaining)"
$ws.Range("D9").Value = "Rajasthan Royals"
$ws.Range("E9").Value = "Sunrisers Hyderabad"
$ws.Range("F9").Value = "Riyan Parag "
$ws.Range("G9").Value = "42"
$ws.Range("H9").Value = "26"
$ws.Range("I9").Value = "2"
$ws.Range("J9").Value = "2"
$ws.Range("K9").Value = "161.53"

